$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.41337446758495844
$ws.Range("B1").Value = 0.41205435375782429
$ws.Range("A2").Value = -0.29411097458942415
$ws.Range("B2").Value = 0.29040247242211592
$ws.Range("A3").Value = -0.18744822605405176
$ws.Range("B3").Value = 0.18621969439606545
$ws.Range("A4").Value = -0.1742196945736918
$ws.Range("B4").Value = 0.17309765719354608
$ws.Range("A5").Value = -0.16709765785917874
$ws.Range("B5").Value = 0.16482943307342435
$ws.Range("A6").Value = -0.10946791955228719
$ws.Range("B6").Value = 0.10928060926347127
$ws.Range("A7").Value = -0.08928061007912369
$ws.Range("B7").Value = 0.088808559704981249
$ws.Range("A8").Value = -0.068808560529618923
$ws.Range("B8").Value = 0.06839987191282848
$ws.Range("A9").Value = -0.062399872623450037
$ws.Range("B9").Value = 0.062051863753144865
$ws.Range("A10").Value = -0.05605186447363053
$ws.Range("B10").Value = 0.056004254364793837
$ws.Range("A11").Value = -0.051504255072600813
$ws.Range("B11").Value = 0.051421532328575381
$ws.Range("A12").Value = -0.045421533052270036
$ws.Range("B12").Value = 0.045160657750289435
$ws.Range("A13").Value = -0.039160658484325594
$ws.Range("B13").Value = 0.039089854533770563
$ws.Range("A14").Value = -0.027089855324580192
$ws.Range("B14").Value = 0.027055913722102254
$ws.Range("A15").Value = -0.021055914461619807
$ws.Range("B15").Value = 0.021029223434584132
$ws.Range("A16").Value = -0.015029224176228206
$ws.Range("B16").Value = 0.015004426336773902
$ws.Range("A17").Value = -0.0090044270813143257
$ws.Range("B17").Value = 0.0089999992280036523
$ws.Range("A18").Value = -0.10794337983039526
$ws.Range("B18").Value = 0.10779366460029038
$ws.Range("A19").Value = -0.027097038229953352
$ws.Range("B19").Value = 0.02701347497041473
$ws.Range("A20").Value = -0.018013475631425635
$ws.Range("B20").Value = 0.018004284326163855
$ws.Range("A21").Value = -0.0090042849879781173
$ws.Range("B21").Value = 0.0089999993376146392
$ws.Range("A22").Value = -0.093954247041962091
$ws.Range("B22").Value = 0.093639429859102563
$ws.Range("A23").Value = -0.084639430538945071
$ws.Range("B23").Value = 0.084127709199825595
$ws.Range("A24").Value = -0.042127710184800549
$ws.Range("B24").Value = 0.041999999009869171
$ws.Range("A25").Value = -0.049361514208630553
$ws.Range("B25").Value = 0.04928083071254008
$ws.Range("A26").Value = -0.043280831401173003
$ws.Range("B26").Value = 0.043180052098691135
$ws.Range("A27").Value = -0.037180052788942763
$ws.Range("B27").Value = 0.036845034918436159
$ws.Range("A28").Value = -0.030845035614800231
$ws.Range("B28").Value = 0.030625943096981523
$ws.Range("A29").Value = -0.018625943850585358
$ws.Range("B29").Value = 0.01853559455513043
$ws.Range("A30").Value = 0.0014644046186296222
$ws.Range("B30").Value = -0.0014748694150199526
$ws.Range("A31").Value = -0.040745076930313928
$ws.Range("B31").Value = 0.040692608080318493
$ws.Range("A32").Value = -0.019692608916892063
$ws.Range("B32").Value = 0.019675579393311438
